# Remove the "ion" (charged species) rows from the molecule table on Sheet1.
# These rows are identified by their row number in the original layout; as
# each row is deleted the rows below it shift up, so we delete from the
# bottom of the list upward to keep the remaining row numbers stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rowsToDelete = @(51,46,45,44,41,40,39,38,35,34,27,23,22,20,19,17,7,4)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$ws.Range("C4").Select()
